$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the "Anio_fundacion" column (column D). Excel's own column
#    delete collapses the now-unused shared string and re-indexes the
#    rest of the table, which is exactly what the target workbook shows
#    ("Anio_fundacion" removed, "Municipio" shifted down to index 33).
# ------------------------------------------------------------------
$ws.Columns.Item(4).Delete()

# ------------------------------------------------------------------
# 2. Widen column A to fit the new, longer municipality codes.
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.25

# ------------------------------------------------------------------
# 3. Append the 14 new corregimiento/municipio rows (33-46).
#    Column A is written first, in the exact order the values were
#    first entered (this controls the order new strings land in the
#    shared-string table: ALTO_BAUDO, BOJAYA, BAJO_BAUDO, ATRATO, ...).
# ------------------------------------------------------------------
$names = @(
  "ALTO_BAUDO",
  "BOJAYA",
  "BAJO_BAUDO",
  "ATRATO",
  "RIO_IRO",
  "BAHIA_SOLANO",
  "MEDIO_BAUDO",
  "EL_LITORAL_DEL_SAN_JUAN",
  "RIO_QUITO",
  "MEDIO_SAN_JUAN",
  "EL_CANTON_DEL_SAN_PABLO",
  "UNION_PANAMERICANA",
  "CARMEN_DEL_DARIEN",
  "MEDIO_ATRATO"
)
$nameRows = @(33, 35, 34, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46)

for ($i = 0; $i -lt $names.Count; $i++) {
  $ws.Cells.Item($nameRows[$i], 1).Value = $names[$i]
}

# Area_km2 (B) / Habitantes_2018 (C) values, row by row.
$data = @(
  @(33, 1532, 1569),
  @(34, 3424, 17402),
  @(35, 3693, 1099),
  @(36, 1017, 13819),
  @(37, 520, 9695),
  @(38, 1667, 9327),
  @(39, 4840, 13560),
  @(40, 3755, 15251),
  @(41, 700, 8961),
  @(42, 620, 15945),
  @(43, 386, 7970),
  @(44, 147, 9592),
  @(45, 4700, 11916),
  @(46, 562, 29489)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
}

# ------------------------------------------------------------------
# 4. Match the final cursor/selection position used when the sheet was
#    saved (bottom of the newly-added data).
# ------------------------------------------------------------------
$ws.Range("B46").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
